$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text values (names, URLs, percentage strings) assign directly -
# Excel keeps these as text since they are not parseable as numbers.
function Set-Text {
    param($cellRef, $text)
    $ws.Range($cellRef).Value = $text
}

# Numeric-looking price strings (e.g. "1.007") must be forced to stay as
# literal text (matching the source data, which is scraped/displayed
# text, not a real number) - otherwise Excel's smart typing would
# silently convert them into floating point numbers and drop formatting
# such as trailing zeros. Setting NumberFormat to "@" (Text) first
# prevents that conversion.
function Set-TextForced {
    param($cellRef, $text)
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
}

# Row 2 - Bitcoin
Set-Text "D2" "27.893.59"
Set-Text "E2" "  +1.49%  "

# Row 3 - Ethereum
Set-Text "D3" "1.906.28"
Set-Text "E3" "  +2.53%  "

# Row 4 - TetherUSD
Set-Text "E4" "  -0.09%  "

# Row 5 - BNB
Set-TextForced "D5" "317.00"
Set-Text "E5" "  +2.03%  "

# Row 6 - USDC
Set-TextForced "D6" "1.007"
Set-Text "E6" "  -0.21%  "

# Row 7 - XRP
Set-TextForced "D7" "0.4835"
Set-Text "E7" "  +1.50%  "

# Row 8 - Cardano
Set-TextForced "D8" "0.3794"
Set-Text "E8" "  +0.05%  "

# Row 9 - Dogecoin
Set-Text "E9" "  +0.61%  "

# Row 10 - Polygon
Set-TextForced "D10" "0.9303"

# Row 11 - Solana
Set-Text "E11" "  +0.30%  "

# Row 12 - TRON
Set-TextForced "D12" "0.07744"
Set-Text "E12" "  -0.42%  "

# Row 13 - WrappedEther
Set-Text "D13" "1.929.49"
Set-Text "E13" "  +3.85%  "

# Row 14 - Polkadot
Set-TextForced "D14" "5.481"
Set-Text "E14" "  +0.74%  "

# Row 15 - Chainlink
Set-TextForced "D15" "6.618"
Set-Text "E15" "  +0.98%  "

# Row 16 - Litecoin
Set-TextForced "D16" "91.78"
Set-Text "E16" "  +1.87%  "

# Row 17 - BinanceUSD
Set-TextForced "D17" "1.009"
Set-Text "E17" "  -0.19%  "

# Row 18 - ShibaInu
Set-TextForced "D18" "0.000008847"
Set-Text "E18" "  +0.38%  "

# Row 19 - Dai
Set-TextForced "D19" "1.007"
Set-Text "E19" "  -0.15%  "

# Row 20 - WrappedBTC
Set-Text "D20" "27.952.66"
Set-Text "E20" "  +1.64%  "

# Row 21 - Avalanche
Set-TextForced "D21" "14.65"
Set-Text "E21" "  +0.12%  "

# Row 22 - Uniswap
Set-TextForced "D22" "5.157"
Set-Text "E22" "  +1.39%  "

# Row 23 - WrappedliquidstakedEther2.0
Set-Text "D23" "2.134.66"
Set-Text "E23" "  +1.76%  "

# Row 24 - Cosmos
Set-TextForced "D24" "10.88"
Set-Text "E24" "  +1.71%  "

# Row 25 - Toncoin
Set-TextForced "D25" "1.919"
Set-Text "E25" "  -0.83%  "

# Row 26 - Monero
Set-TextForced "D26" "154.74"
Set-Text "E26" "  -0.49%  "

# Row 27 - EthereumClassic
Set-TextForced "D27" "18.46"
Set-Text "E27" "  -0.01%  "

# Row 28 - LidoDAOToken
Set-TextForced "D28" "2.127"
Set-Text "E28" "  +6.13%  "

# Row 29 - BitcoinCash
Set-TextForced "D29" "117.23"
Set-Text "E29" "  +1.69%  "

# Row 30 - InternetComputer(DFINITY)
Set-TextForced "D30" "4.954"
Set-Text "E30" "  +0.18%  "

# Row 31 - Stellar
Set-Text "E31" "  +1.19%  "

# Row 32 - HuobiToken
Set-TextForced "D32" "3.237"
Set-Text "E32" "  -2.76%  "

# Row 33 - ARBITRUM
Set-TextForced "D33" "1.252"
Set-Text "E33" "  +4.26%  "

# Row 34 - ImmutableX
Set-TextForced "D34" "0.7668"
Set-Text "E34" "  +2.23%  "

# Row 35 - Filecoin
Set-TextForced "D35" "4.657"
Set-Text "E35" "  +1.72%  "

# Row 36 - VeChain
Set-TextForced "D36" "0.02045"
Set-Text "E36" "  +0.10%  "

# Row 37 - RenderToken
Set-TextForced "D37" "2.529"
Set-Text "E37" "  -6.49%  "

# Row 38 - TrustWalletToken
Set-TextForced "D38" "1.097"
Set-Text "E38" "  -1.96%  "

# Row 39 - was Hedera, now MXToken (rows 39/40 content swapped)
Set-Text "B39" "MXToken"
Set-Text "C39" "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextForced "D39" "3.002"
Set-Text "E39" "  +0.78%  "

# Row 40 - was MXToken, now Hedera
Set-Text "B40" "Hedera"
Set-Text "C40" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextForced "D40" "0.05276"
Set-Text "E40" "  -0.25%  "

# Row 41 - TheSandbox
Set-TextForced "D41" "0.5469"
Set-Text "E41" "  -1.32%  "

# Row 42 - FraxShare
Set-TextForced "D42" "6.953"
Set-Text "E42" "  -1.04%  "

# Row 43 - Algorand
Set-TextForced "D43" "0.1525"
Set-Text "E43" "  +0.42%  "

# Row 44 - Aptos
Set-TextForced "D44" "8.402"
Set-Text "E44" "  -1.58%  "

# Row 45 - was Quant, now EnergySwap (rows 45/46 content swapped)
Set-Text "B45" "EnergySwap"
Set-Text "C45" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextForced "D45" "10.69"
Set-Text "E45" "  +0.11%  "

# Row 46 - was EnergySwap, now Quant
Set-Text "B46" "Quant"
Set-Text "C46" "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextForced "D46" "108.96"
Set-Text "E46" "  +4.88%  "

# Row 47 - Decentraland
Set-TextForced "D47" "0.4804"
Set-Text "E47" "  -1.11%  "

# Row 48 - PaxDollar
Set-Text "E48" "  -0.19%  "

# Row 49 - NEARProtocol
Set-TextForced "D49" "1.650"
Set-Text "E49" "  -0.82%  "

# Row 50 - Aave
Set-TextForced "D50" "67.68"
Set-Text "E50" "  +0.57%  "

# Row 51 - Cronos
Set-TextForced "D51" "0.06090"
Set-Text "E51" "  -0.15%  "
